$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source data stores every "Price" cell as text (even values that look
# like plain numbers, e.g. "1.00"), so force text format on each Price cell
# before writing its new value -- this avoids Excel auto-converting strings
# such as "1.00" into the number 1 and losing the trailing zero.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '68.210.14'
$ws.Range("E2").Value = '  -3.82%  '

# Row 3
$ws.Range("D3").Value = '3.691.37'
$ws.Range("E3").Value = '  -4.11%  '

# Row 4
$ws.Range("E4").Value = '  -0.22%  '

# Row 5
$ws.Range("D5").Value = '600.44'
$ws.Range("E5").Value = '  +1.39%  '

# Row 6
$ws.Range("D6").Value = '181.73'
$ws.Range("E6").Value = '  +8.84%  '

# Row 7
$ws.Range("D7").Value = '3.683.57'
$ws.Range("E7").Value = '  -4.25%  '

# Row 8
$ws.Range("D8").Value = '0.627'
$ws.Range("E8").Value = '  -6.50%  '

# Row 9
$ws.Range("D9").Value = '1.00'
$ws.Range("E9").Value = '  +0.11%  '

# Row 10
$ws.Range("D10").Value = '0.718'
$ws.Range("E10").Value = '  -3.99%  '

# Row 11
$ws.Range("D11").Value = '0.162'
$ws.Range("E11").Value = '  -6.40%  '

# Row 12
$ws.Range("D12").Value = '55.53'
$ws.Range("E12").Value = '  +4.38%  '

# Row 13
$ws.Range("D13").Value = '0.0000289'
$ws.Range("E13").Value = '  -9.43%  '

# Row 14
$ws.Range("D14").Value = '10.41'
$ws.Range("E14").Value = '  -6.65%  '

# Row 15
$ws.Range("D15").Value = '4.282.18'
$ws.Range("E15").Value = '  -4.45%  '

# Row 16
$ws.Range("D16").Value = '3.690.28'
$ws.Range("E16").Value = '  -4.88%  '

# Row 17
$ws.Range("D17").Value = '19.36'
$ws.Range("E17").Value = '  -6.14%  '

# Row 18
$ws.Range("E18").Value = '  -2.30%  '

# Row 19
$ws.Range("D19").Value = '1.12'
$ws.Range("E19").Value = '  -6.44%  '

# Row 20
$ws.Range("D20").Value = '12.79'
$ws.Range("E20").Value = '  -7.04%  '

# Row 21
$ws.Range("D21").Value = '67.937.56'
$ws.Range("E21").Value = '  -4.14%  '

# Row 22
$ws.Range("D22").Value = '409.09'
$ws.Range("E22").Value = '  -5.45%  '

# Row 23
$ws.Range("D23").Value = '4.57'
$ws.Range("E23").Value = '  -2.96%  '

# Row 24
$ws.Range("D24").Value = '88.24'
$ws.Range("E24").Value = '  -6.21%  '

# Row 25
$ws.Range("D25").Value = '3.01'
$ws.Range("E25").Value = '  -7.78%  '

# Row 26
$ws.Range("D26").Value = '12.76'
$ws.Range("E26").Value = '  -7.05%  '

# Row 27
$ws.Range("D27").Value = '10.96'
$ws.Range("E27").Value = '  +0.66%  '

# Row 28
$ws.Range("D28").Value = '3.85'
$ws.Range("E28").Value = '  -5.98%  '

# Row 29
$ws.Range("D29").Value = '6.06'
$ws.Range("E29").Value = '  +2.25%  '

# Row 30
$ws.Range("D30").Value = '9.47'
$ws.Range("E30").Value = '  -6.88%  '

# Row 31
$ws.Range("D31").Value = '32.69'
$ws.Range("E31").Value = '  -6.33%  '

# Row 32
$ws.Range("D32").Value = '7.37'
$ws.Range("E32").Value = '  -7.07%  '

# Row 33
$ws.Range("D33").Value = '12.46'
$ws.Range("E33").Value = '  -7.49%  '

# Row 34
$ws.Range("E34").Value = '  -6.20%  '

# Row 35
$ws.Range("D35").Value = '43.77'
$ws.Range("E35").Value = '  -10.82%  '

# Row 36
$ws.Range("D36").Value = '64.72'
$ws.Range("E36").Value = '  -5.59%  '

# Row 37
$ws.Range("D37").Value = '589.88'
$ws.Range("E37").Value = '  -4.55%  '

# Row 38
$ws.Range("D38").Value = '0.0₃0880'
$ws.Range("E38").Value = '  -9.58%  '

# Row 39
$ws.Range("B39").Value = 'Dai'
$ws.Range("C39").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D39").Value = '0.999'
$ws.Range("E39").Value = '  +0.01%  '

# Row 40
$ws.Range("B40").Value = 'TheGraph'
$ws.Range("C40").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D40").Value = '0.399'
$ws.Range("E40").Value = '  -4.41%  '

# Row 41
$ws.Range("D41").Value = '1.00'
$ws.Range("E41").Value = '  +0.32%  '

# Row 42
$ws.Range("E42").Value = '  -3.98%  '

# Row 43
$ws.Range("B43").Value = 'Fetch.AI'
$ws.Range("C43").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D43").Value = '2.71'
$ws.Range("E43").Value = '  +2.72%  '

# Row 44
$ws.Range("B44").Value = 'ThetaToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D44").Value = '2.99'
$ws.Range("E44").Value = '  -8.53%  '

# Row 45
$ws.Range("D45").Value = '0.0434'
$ws.Range("E45").Value = '  -6.75%  '

# Row 46
$ws.Range("D46").Value = '2.89'
$ws.Range("E46").Value = '  -12.50%  '

# Row 47
$ws.Range("D47").Value = '9.28'
$ws.Range("E47").Value = '  -7.92%  '

# Row 48
$ws.Range("D48").Value = '2.799.09'
$ws.Range("E48").Value = '  -0.79%  '

# Row 49
$ws.Range("E49").Value = '  -6.74%  '

# Row 50
$ws.Range("E50").Value = '  -3.21%  '

# Row 51
$ws.Range("D51").Value = '3.13'
$ws.Range("E51").Value = '  -5.48%  '
